$wb = $excel.ActiveWorkbook

# Rename the first two sheets.
$sheet0 = $wb.Worksheets.Item(1)
$sheet1 = $wb.Worksheets.Item(2)
$sheet0.Name = "sheet0"
$sheet1.Name = "sheet1"

# Duplicate sheet0 -> sheet0_bad (appended at the end), then shift its
# contents down one row and right one column so column/row 1 is blank --
# a "bad" copy used to test the "first column" option.
$sheet0.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet0Bad = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet0Bad.Name = "sheet0_bad"
$sheet0Bad.Rows(1).Insert()
$sheet0Bad.Columns(1).Insert()

# Duplicate sheet1 -> sheet1_bad the same way.
$sheet1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet1Bad = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet1Bad.Name = "sheet1_bad"
$sheet1Bad.Rows(1).Insert()
$sheet1Bad.Columns(1).Insert()

# Selections on the new sheets, matching what a user left behind.
$sheet0Bad.Range("E30").Select()
$sheet1Bad.Range("G18").Select()

# sheet1_bad ends up the active / selected tab.
$sheet1Bad.Activate()
